$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("D2").Value = 909
$ws.Range("E2").Value = 77
$ws.Range("F2").Value = 77
$ws.Range("G2").Value = 67
$ws.Range("H2").Value = 38
$ws.Range("I2").Value = 38
$ws.Range("J2").ClearContents()
$ws.Range("K2").Value = 1390
$ws.Range("L2").Value = 687
$ws.Range("M2").Value = 702
$ws.Range("N2").Value = 702
$ws.Range("O2").ClearContents()
$ws.Range("P2").Value = 113
$ws.Range("Q2").Value = 123
$ws.Range("R2").Value = -49
$ws.Range("S2").Value = -35
$ws.Range("T2").Value = 53
$ws.Range("U2").Value = 71
$ws.Range("V2").Value = 355
$ws.Range("W2").Value = 8.43
$ws.Range("X2").Value = 4.21
$ws.Range("Y2").Value = 5.47
$ws.Range("Z2").Value = 2.8
$ws.Range("AA2").Value = 97.86
$ws.Range("AB2").Value = 696.33
$ws.Range("AC2").Value = 327
$ws.Range("AD2").Value = 22.06
$ws.Range("AE2").Value = 6007
$ws.Range("AF2").Value = 1.2
$ws.Range("AG2").Value = 96
$ws.Range("AH2").Value = 1.34
$ws.Range("AI2").Value = 29.49
$ws.Range("AJ2").Value = 11696733

# --- Row 3 ---
$ws.Range("D3").Value = 1096
$ws.Range("E3").Value = 142
$ws.Range("F3").Value = 142
$ws.Range("G3").Value = 132
$ws.Range("H3").Value = 108
$ws.Range("I3").Value = 108
$ws.Range("J3").ClearContents()
$ws.Range("K3").Value = 1410
$ws.Range("L3").Value = 604
$ws.Range("M3").Value = 807
$ws.Range("N3").Value = 807
$ws.Range("O3").ClearContents()
$ws.Range("P3").Value = 113
$ws.Range("Q3").Value = 88
$ws.Range("R3").Value = -42
$ws.Range("S3").Value = -87
$ws.Range("T3").Value = 39
$ws.Range("U3").Value = 48
$ws.Range("V3").Value = 273
$ws.Range("W3").Value = 13
$ws.Range("X3").Value = 9.81
$ws.Range("Y3").Value = 14.25
$ws.Range("Z3").Value = 7.68
$ws.Range("AA3").Value = 74.84
$ws.Range("AB3").Value = 787.8
$ws.Range("AC3").Value = 919
$ws.Range("AD3").Value = 9.01
$ws.Range("AE3").Value = 6899
$ws.Range("AF3").Value = 1.2
$ws.Range("AG3").Value = 145
$ws.Range("AH3").Value = 1.75
$ws.Range("AI3").Value = 15.73
$ws.Range("AJ3").Value = 11696733

# --- Row 4 ---
$ws.Range("D4").Value = 1084
$ws.Range("E4").Value = 64
$ws.Range("F4").Value = 64
$ws.Range("G4").Value = 59
$ws.Range("H4").Value = 47
$ws.Range("I4").Value = 47
$ws.Range("J4").ClearContents()
$ws.Range("K4").Value = 1465
$ws.Range("L4").Value = 617
$ws.Range("M4").Value = 848
$ws.Range("N4").Value = 848
$ws.Range("O4").ClearContents()
$ws.Range("P4").Value = 113
$ws.Range("Q4").Value = 118
$ws.Range("R4").Value = -35
$ws.Range("S4").Value = -1
$ws.Range("T4").Value = 39
$ws.Range("U4").Value = 78
$ws.Range("V4").Value = 287
$ws.Range("W4").Value = 5.9
$ws.Range("X4").Value = 4.3
$ws.Range("Y4").Value = 5.63
$ws.Range("Z4").Value = 3.24
$ws.Range("AA4").Value = 72.70999999999999
$ws.Range("AB4").Value = 824.85
$ws.Range("AC4").Value = 398
$ws.Range("AD4").Value = 21.59
$ws.Range("AE4").Value = 7254
$ws.Range("AF4").Value = 1.19
$ws.Range("AG4").Value = 96
$ws.Range("AH4").Value = 1.12
$ws.Range("AI4").Value = 24.2
$ws.Range("AJ4").Value = 11696733

# --- Row 5 ---
$ws.Range("D5").Value = 1189
$ws.Range("E5").Value = 147
$ws.Range("F5").Value = 147
$ws.Range("G5").Value = 130
$ws.Range("H5").Value = 74
$ws.Range("I5").Value = 76
$ws.Range("J5").ClearContents()
$ws.Range("K5").Value = 1569
$ws.Range("L5").Value = 589
$ws.Range("M5").Value = 981
$ws.Range("N5").Value = 966
$ws.Range("O5").ClearContents()
$ws.Range("P5").Value = 113
$ws.Range("Q5").Value = 146
$ws.Range("R5").Value = -242
$ws.Range("S5").Value = -26
$ws.Range("T5").Value = 103
$ws.Range("U5").Value = 43
$ws.Range("V5").Value = 215
$ws.Range("W5").Value = 12.33
$ws.Range("X5").Value = 6.19
$ws.Range("Y5").Value = 8.369999999999999
$ws.Range("Z5").Value = 4.85
$ws.Range("AA5").Value = 60.04
$ws.Range("AB5").Value = 941.54
$ws.Range("AC5").Value = 649
$ws.Range("AD5").Value = 30.47
$ws.Range("AE5").Value = 8264
$ws.Range("AF5").Value = 2.39
$ws.Range("AG5").Value = 241
$ws.Range("AH5").Value = 1.22
$ws.Range("AI5").Value = 37.1
$ws.Range("AJ5").Value = 11696733

# --- Row 6 (no J6/O6 originally) ---
$ws.Range("D6").Value = 1340
$ws.Range("E6").Value = 152
$ws.Range("F6").Value = 152
$ws.Range("G6").Value = 152
$ws.Range("H6").Value = 77
$ws.Range("I6").Value = 76
$ws.Range("K6").Value = 1599
$ws.Range("L6").Value = 566
$ws.Range("M6").Value = 1033
$ws.Range("N6").Value = 1018
$ws.Range("P6").Value = 113
$ws.Range("Q6").Value = 144
$ws.Range("R6").Value = -41
$ws.Range("S6").Value = -21
$ws.Range("T6").Value = 32
$ws.Range("U6").Value = 111
$ws.Range("V6").Value = 223
$ws.Range("W6").Value = 11.33
$ws.Range("X6").Value = 5.74
$ws.Range("Y6").Value = 7.71
$ws.Range("Z6").Value = 4.85
$ws.Range("AA6").Value = 54.86
$ws.Range("AB6").Value = 990.61
$ws.Range("AC6").Value = 654
$ws.Range("AD6").Value = 41.56
$ws.Range("AE6").Value = 8705
$ws.Range("AF6").Value = 3.12
$ws.Range("AG6").Value = 241
$ws.Range("AH6").Value = 0.89
$ws.Range("AI6").Value = 36.85
$ws.Range("AJ6").Value = 11696733

# --- Rows 7, 8, 9: remove all data columns D..AJ, keep only A,B,C ---
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
